$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.29903
$ws.Range("H2").Value = 6.89709
$ws.Range("I2").Value = 0.04075801785348079
$ws.Range("J2").Value = 0.04075801785348079
$ws.Range("M2").Value = 7.106976666666665
$ws.Range("N2").Value = 21.32093
$ws.Range("O2").Value = 0.1598176868560746
$ws.Range("P2").Value = 0.1598176868560746
$ws.Range("Q2").Value = 16.33915256596666
$ws.Range("R2").Value = 147.0523730937
$ws.Range("S2").Value = 0.006513852134181889
$ws.Range("T2").Value = 0.006513852134181889

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.29903
$ws.Range("H3").Value = 6.89709
$ws.Range("I3").Value = 0.04075801785348079
$ws.Range("J3").Value = 0.04075801785348079
$ws.Range("O3").Value = 0.6350325402576649
$ws.Range("P3").Value = 0.6350325402576648
$ws.Range("Q3").Value = 64.92331207976667
$ws.Range("R3").Value = 584.3098087179001
$ws.Range("S3").Value = 0.02588266761336316
$ws.Range("T3").Value = 0.02588266761336316

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.29903
$ws.Range("H4").Value = 6.89709
$ws.Range("I4").Value = 0.04075801785348079
$ws.Range("J4").Value = 0.04075801785348079
$ws.Range("O4").Value = 0.2051497728862606
$ws.Range("P4").Value = 0.2051497728862606
$ws.Range("Q4").Value = 20.97373265751667
$ws.Range("R4").Value = 188.76359391765
$ws.Range("S4").Value = 0.00836149810593574
$ws.Range("T4").Value = 0.00836149810593574

# Row 5
$ws.Range("I5").Value = 0.3949230674234065
$ws.Range("J5").Value = 0.3949230674234066
$ws.Range("M5").Value = 7.106976666666665
$ws.Range("N5").Value = 21.32093
$ws.Range("O5").Value = 0.1598176868560746
$ws.Range("P5").Value = 0.1598176868560746
$ws.Range("Q5").Value = 158.3175186204378
$ws.Range("R5").Value = 1424.85766758394
$ws.Range("S5").Value = 0.06311569112171442
$ws.Range("T5").Value = 0.06311569112171442

# Row 6
$ws.Range("I6").Value = 0.3949230674234065
$ws.Range("J6").Value = 0.3949230674234066
$ws.Range("O6").Value = 0.6350325402576649
$ws.Range("P6").Value = 0.6350325402576648
$ws.Range("S6").Value = 0.2507889987122349
$ws.Range("T6").Value = 0.2507889987122349

# Row 7
$ws.Range("I7").Value = 0.3949230674234065
$ws.Range("J7").Value = 0.3949230674234066
$ws.Range("O7").Value = 0.2051497728862606
$ws.Range("P7").Value = 0.2051497728862606
$ws.Range("S7").Value = 0.08101837758945725
$ws.Range("T7").Value = 0.08101837758945725

# Row 8
$ws.Range("I8").Value = 0.5643189147231126
$ws.Range("J8").Value = 0.5643189147231126
$ws.Range("M8").Value = 7.106976666666665
$ws.Range("N8").Value = 21.32093
$ws.Range("O8").Value = 0.1598176868560746
$ws.Range("P8").Value = 0.1598176868560746
$ws.Range("Q8").Value = 226.2252515975633
$ws.Range("R8").Value = 2036.02726437807
$ws.Range("S8").Value = 0.09018814360017827
$ws.Range("T8").Value = 0.09018814360017825

# Row 9
$ws.Range("I9").Value = 0.5643189147231126
$ws.Range("J9").Value = 0.5643189147231126
$ws.Range("O9").Value = 0.6350325402576649
$ws.Range("P9").Value = 0.6350325402576648
$ws.Range("S9").Value = 0.3583608739320667
$ws.Range("T9").Value = 0.3583608739320667

# Row 10
$ws.Range("I10").Value = 0.5643189147231126
$ws.Range("J10").Value = 0.5643189147231126
$ws.Range("O10").Value = 0.2051497728862606
$ws.Range("P10").Value = 0.2051497728862606
$ws.Range("S10").Value = 0.1157698971908676
$ws.Range("T10").Value = 0.1157698971908676
